# Ajout de donnee mise a jour estimation
#
# Applies the commit "Ajout de donnee mise a jour estimation": fills in
# several rows of the "Iteration #3" and "Iteration #4" sheets with dates,
# task descriptions and time estimates, and updates the active-sheet /
# selection state so that "Iteration #4" ends up the active tab.

$wb = $excel.ActiveWorkbook

$wsEstimation = $wb.Worksheets.Item("Estimation")
$wsIter1      = $wb.Worksheets.Item("Iteration #1")
$wsIter2      = $wb.Worksheets.Item("Iteration #2")
$wsIter3      = $wb.Worksheets.Item("Iteration #3")
$wsIter4      = $wb.Worksheets.Item("Iteration #4")

# ---------------------------------------------------------------------
# Iteration #3 : rows 14-24 (dates use the Excel 1900 date-system serial
# numbers taken from the source diff so they match exactly).
# ---------------------------------------------------------------------
$wsIter3.Range("A14").Value = 43213
$wsIter3.Range("B14").Value = "Travailler pour bien faire marcher le gyroscope"
$wsIter3.Range("C14").Value = "3h30"

$wsIter3.Range("A15").Value = 43214
$wsIter3.Range("B15").Value = "Arrêt du projet de la boule"
$wsIter3.Range("C15").Value = "1h00"

$wsIter3.Range("A16").Value = 43214
$wsIter3.Range("B16").Value = "Debut projet de musique"
$wsIter3.Range("C16").Value = "2h00"

# row 17 is left untouched (stays blank)

$wsIter3.Range("A18").Value = 43219
$wsIter3.Range("B18").Value = 'Création d''image pause et play en paint '
$wsIter3.Range("C18").Value = "30 min"

$wsIter3.Range("A19").Value = 43220
$wsIter3.Range("B19").Value = "Faire marcher bouton play,pause et stop"
$wsIter3.Range("C19").Value = "3h00"

$wsIter3.Range("A20").Value = 43191
$wsIter3.Range("B20").Value = "Modifier code pour ajouter une autre musique et finir faire fonctionner le bouton stop"
$wsIter3.Range("C20").Value = "3h00"

$wsIter3.Range("A21").Value = 43222
$wsIter3.Range("B21").Value = "Créer image avancer,reculer,passer et avancer."
$wsIter3.Range("C21").Value = "1h00"

# row 22 : A22 stays blank, only B22/C22 are filled in
$wsIter3.Range("B22").Value = "Faire marcher les boutons avancer,reculer,passer et avancer"
$wsIter3.Range("C22").Value = "4h00"

$wsIter3.Range("A23").Value = 43225
$wsIter3.Range("B23").Value = 'Penser a quoi rajouter,commencer créer un bouton ajouter musique.Tester des trucs '
$wsIter3.Range("C23").Value = "3h00"

# row 24 : A24 and C24 stay blank, only B24 is filled in
$wsIter3.Range("B24").Value = "dans la classe"

# ---------------------------------------------------------------------
# Iteration #4 : row 14
# ---------------------------------------------------------------------
$wsIter4.Range("A14").Value = 43226
$wsIter4.Range("B14").Value = "Penser comment faire marcher la classe ajoutermusique"
$wsIter4.Range("C14").Value = "1h00"

# ---------------------------------------------------------------------
# View / selection state. Selecting a range activates its sheet, so the
# sheets are touched in the same order as the target workbook (ending on
# "Iteration #4" so it becomes the active / tabSelected sheet, matching
# activeTab="4" on the workbook and tabSelected="1" on that sheet).
# ---------------------------------------------------------------------
$wsEstimation.Activate()
$wsEstimation.Range("D24").Select()

$wsIter1.Activate()
$wsIter1.Range("B21").Select()

$wsIter2.Activate()
$wsIter2.Range("C19").Select()

$wsIter3.Activate()
$wsIter3.Range("B25").Select()

$wsIter4.Activate()
$wsIter4.Range("B14").Select()
